# Defect Log.xlsx update
# Nhom ThanhChV, HuyDV, DucNH, LinhTA nop bai Defect Log.xlsx
#
# Fills in rows 15-25 (Created Date / Title / Description / Status) on the
# active sheet and moves the active selection to K16, matching the target
# workbook. Shared-string cells are written in the exact order the new
# strings are first introduced so that the resulting shared strings table
# matches the canonical ordering (indices 121-139).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C / D text values, written in "first use" order so the
#     resulting shared-strings table lines up index-for-index with the
#     target workbook (121 .. 139). Reused strings (same Title repeated on
#     consecutive rows) are simply referenced again later below.

$ws.Range("C15").Value = "Language manager "
$ws.Range("D15").Value = "Dư column Version "
$ws.Range("D16").Value = "Dư column Date"
$ws.Range("C17").Value = "New & Edit user"
$ws.Range("D17").Value = "Dư phần Basic Setting  "
$ws.Range("D18").Value = "Dư các  group user :  publisher, editor, super administrator, author, editor."
$ws.Range("C19").Value = "New & Edit user(mapping items to db)"
$ws.Range("D19").Value = "Dư column Cofirm password trong phần mapping to db"
$ws.Range("C20").Value = "new & edit user"
$ws.Range("D20").Value = "Dùng các radion button thay cho cac check box trong Assigned User Groups"
$ws.Range("C21").Value = "Contend manager _ Article manager(Item Definition)"
$ws.Range("C22").Value = "Contend manager_ Article manager(events)"
$ws.Range("D22").Value = "Mô tả sai event của btnpublishArticle"
$ws.Range("D23").Value = "Mô tả sai event của btnunpublishArticle"
$ws.Range("C24").Value = "Contend manager_ Article manager(item definition)"
$ws.Range("D24").Value = "các items : pageNumberList, itemNotes có type là : html, mô tả type của item  chưa rõ.!"
$ws.Range("C25").Value = "Contend manager_ Article manager(event)"
$ws.Range("D25").Value = "Mô tả chưa chính xác event của btnreset"
$ws.Range("D21").Value = "mô tả thiếu item btnhelp trong phần item definiton"

# --- Reused Title values (repeat strings already introduced above)

$ws.Range("C16").Value = "Language manager "
$ws.Range("C18").Value = "New & Edit user"
$ws.Range("C23").Value = "Contend manager_ Article manager(events)"

# --- Status column (reuses existing shared string "Error") and the
#     Created Date column for every newly-populated row.

$rows = 15..25
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "Error"
    $ws.Range("B$r").Value = 40837
}

# --- Restore the active cell selection recorded in the saved workbook.

$ws.Range("K16").Select()
